# Apply updated employee absence data values to rows 2-11
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 82237; B = "Heloísa Marques";          C = "P&D";                     D = "Problemas pessoais"; E = 8; F = 45079; G = 9770.16 },
    @{ Row = 3;  A = 60457; B = "Matteo Rocha";              C = "Financeiro";               D = "Problemas pessoais"; E = 1; F = 45104; G = 7171.14 },
    @{ Row = 4;  A = 8835;  B = "Josué Cirino";              C = "Atendimento ao Cliente";   D = "Viagem de negocios"; E = 1; F = 45082; G = 5994.21 },
    @{ Row = 5;  A = 91157; B = "Ana Julia Cavalcante";      C = "Operacoes";                D = "Viagem de negocios"; E = 7; F = 45081; G = 6202.25 },
    @{ Row = 6;  A = 97591; B = "Sr. Yago Araújo";           C = "Operacoes";                D = "Consulta medica";    E = 5; F = 45096; G = 5530.53 },
    @{ Row = 7;  A = 73184; B = "Maya da Cunha";             C = "Atendimento ao Cliente";   D = "Consulta medica";    E = 4; F = 45102; G = 6779.41 },
    @{ Row = 8;  A = 86109; B = "Clarice Marques";           C = "Engenharia";               D = "Consulta medica";    E = 8; F = 45102; G = 7273.03 },
    @{ Row = 9;  A = 80625; B = "Bárbara Caldeira";          C = "Juridico";                 D = "Outros";             E = 6; F = 45092; G = 4408.59 },
    @{ Row = 10; A = 81353; B = "Lucas Gabriel Gonçalves";   C = "Marketing";                D = "Doenca";             E = 2; F = 45092; G = 4973.45 },
    @{ Row = 11; A = 25260; B = "Paulo da Cunha";            C = "Operacoes";                D = "Consulta medica";    E = 8; F = 45096; G = 8040.68 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
